# Fix HPO code formatting: "HP.xxxxxxx" -> "HP:xxxxxxx"
# on the "Include from hp.owl" worksheet, column A (HPO concept codes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Include from hp.owl")

$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val -like "HP.*") {
        $cell.Value = $val -replace "^HP\.", "HP:"
    }
}
